# Weekly update: prepend a new week's price block (date 45075) for
# "Agrícola del Norte S.A. de Arica - Pimiento" and push the existing
# data down by 6 rows (one full Primera/Segunda/Tercera x rojo/verde block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows right before the first data row of the block (row 977),
# shifting all existing rows (977:1006) down to (983:1012).
$ws.Rows("977:982").Insert()

# Shared constant values across every data row in this sheet.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$catId     = 100112002
$categoria = "Pimiento"
$unidad    = "$/caja 15 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidades = 15
$clasificacion = "Hortaliza"
$fecha = 45075

# New week's data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Precio$/Kg
$newRows = @(
    @{ Row = 977; Variedad = "Zafiro rojo";  Calidad = "Primera"; J = 340; K = 15000; L = 16000; M = 15588; P = 1039 },
    @{ Row = 978; Variedad = "Zafiro rojo";  Calidad = "Segunda"; J = 250; K = 13000; L = 14000; M = 13520; P = 901 },
    @{ Row = 979; Variedad = "Zafiro rojo";  Calidad = "Tercera"; J = 220; K = 11000; L = 12000; M = 11409; P = 761 },
    @{ Row = 980; Variedad = "Zafiro verde"; Calidad = "Primera"; J = 140; K = 9000;  L = 10000; M = 9571;  P = 638 },
    @{ Row = 981; Variedad = "Zafiro verde"; Calidad = "Segunda"; J = 130; K = 7000;  L = 8000;  M = 7538;  P = 503 },
    @{ Row = 982; Variedad = "Zafiro verde"; Calidad = "Tercera"; J = 120; K = 5000;  L = 6000;  M = 5500;  P = 367 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $catId
    $ws.Cells.Item($row, 7).Value2  = $categoria
    $ws.Cells.Item($row, 8).Value2  = $r.Variedad
    $ws.Cells.Item($row, 9).Value2  = $r.Calidad
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $unidad
    $ws.Cells.Item($row, 15).Value2 = $origen
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $kgUnidades
    $ws.Cells.Item($row, 18).Value2 = $clasificacion
}
